$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("B2").Value = 70.06
$ws.Range("C2").Value = "The applicant has strong skills in ReactJS and JavaScript, which are required for the job. Additionally, the applicant has experience in web development using Django and Python, which can be beneficial for the role. However, the job requires expertise in MongoDB and NodeJS, which are not explicitly mentioned in the applicant's projects or skills. Furthermore, the applicant lacks specific experience in MongoDB and NodeJS, which are essential for the position. Although the applicant's projects showcase proficiency in ReactJS and JavaScript, the absence of direct experience in MongoDB and NodeJS could impact the applicant's suitability for the role."
$ws.Range("E2").Value = "Based on the applicant's responses and demonstrated adaptability, it is evident that the applicant possesses strong communication and problem-solving skills. The applicant's willingness to learn from the Japanese work culture and the ability to adapt to diverse cultures indicate a high level of adaptability and open-mindedness. Furthermore, the applicant's conscious effort to seek feedback and implement it, as well as the desire for continuous improvement, reflects a strong commitment to personal and professional growth. Overall, the applicant demonstrates a willingness to work in Japan and exhibits positive personality traits, warranting a score of 4 out of 5."

# Update row 3
$ws.Range("B3").Value = 90.53
$ws.Range("C3").Value = "The applicant has a strong set of technical skills including ReactJS, NodeJS, ExpressJS, HTML, CSS, and WebRTC which are relevant to the job description. However, the applicant lacks experience with MongoDB and could benefit from gaining more exposure to this technology. Although the projects the applicant has worked on demonstrate proficiency in web and mobile application development, the absence of MongoDB in the projects may raise some concerns about the candidate's hands-on experience with the technology required by the company. Overall, the applicant's technical skills align well with the job requirements, but there is room for improvement in MongoDB."
$ws.Range("E3").Value = "The applicant has expressed a strong interest in working in Japan, highlighting the appeal of Japan's work culture and expressing a willingness to adapt to a new environment. The applicant has also demonstrated the ability to work collaboratively and has highlighted the value of teamwork and open communication. These qualities align with the soft skills the company is seeking, indicating a good fit in terms of adaptability and teamwork. Overall, the applicant's responses and attitudes suggest a positive willingness to work in Japan, earning a score of 4."

# Update row 4
$ws.Range("B4").Value = 76.82
$ws.Range("C4").Value = "The applicant possesses strong technical skills in Python, Django Rest Framework, ReactJS, PyTorch, and Tensorflow, which demonstrate proficiency in web and machine learning technologies. However, the job requires expertise in MongoDB, ExpressJS, and NodeJS, which are not evident from the applicant's projects. Additionally, the applicant lacks experience in web/mobile application development using the MERN stack as per the company's requirements."
$ws.Range("E4").Value = "The applicant's responses reflect a high level of adaptability, as well as a positive attitude towards learning and working within a new cultural setting like Japan. The applicant's strong teamwork skills and willingness to address weaknesses indicate a good fit for the company's requirements. Overall, the applicant's positive attitude and adaptability make them suitable for the role and working in Japan."

# Add new row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 75.82
$ws.Range("C5").Value = "The applicant has strong skills in JavaScript, HTML, CSS, and API integration, which are essential for web development. While the applicant's projects do not specifically mention MERN stack technologies, the experience in developing web and mobile applications, messaging, file sharing, and screen sharing/recording demonstrate a strong foundation in web development. The applicant's background in engineering and graphic design sets them apart and showcases their ability to bridge technical and creative aspects in projects."
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = "The applicant has demonstrated strong communication skills, adaptability, and a willingness to work in Japan. They have expressed a keen interest in Japan's AI field and work culture. Additionally, the applicant has emphasized their ability to adapt to new cultures and learn a new language while working in Japan. Their commitment to collaborating within a team and addressing weaknesses also shows a positive attitude towards personal and professional growth."
